# Insert a new data row at row 650 (Femacal de La Calera - Ají sheet).
# This shifts the existing rows 650-730 down to 651-731 and grows the
# sheet's used range from A1:R730 to A1:R731, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("650:650").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A650").Value = 3
$ws.Range("B650").Value = 'Femacal de La Calera'
$ws.Range("C650").Value = 'Coquimbo'
$ws.Range("D650").Value = 45124
$ws.Range("E650").Value = 5
$ws.Range("F650").Value = 100112021
$ws.Range("G650").Value = 'Ají'
$ws.Range("H650").Value = 'Inferno'
$ws.Range("I650").Value = 'Primera'
$ws.Range("J650").Value = 83
$ws.Range("K650").Value = 11000
$ws.Range("L650").Value = 11500
$ws.Range("M650").Value = 11229
$ws.Range("N650").Value = '$/caja 10 kilos'
$ws.Range("O650").Value = 'Región de Arica y Parinacota'
$ws.Range("P650").Value = 1123
$ws.Range("Q650").Value = 10
$ws.Range("R650").Value = 'Hortaliza'
